$d = $word.ActiveDocument

# --- Step 1: insert a "Purpose: " run before the existing "Consistency of nomenclature..." run ---
$full = $d.Content.Text
$idx = $full.IndexOf("Consistency of nomenclature")
$r = $d.Range($idx, $idx)
$r.InsertBefore("Purpose: ")

# --- Step 2: split the paragraph just before " Recommendations for standardized..." and drop the
#     leading space, starting a new paragraph that begins with "Recommendations" ---
$full = $d.Content.Text
$idx = $full.IndexOf(" Recommendations for standardized nomenclature")
$spaceRange = $d.Range($idx, $idx + 1)
$spaceRange.Text = ""

$full = $d.Content.Text
$idx = $full.IndexOf("Recommendations for standardized nomenclature")
$splitPoint = $d.Range($idx, $idx)
$splitPoint.InsertParagraphBefore()

# --- Step 3: add the "Methods and Materials: " label run before "Recommendations" ---
$full = $d.Content.Text
$idx = $full.IndexOf("Recommendations for standardized nomenclature")
$r = $d.Range($idx, $idx)
$r.InsertBefore("Methods and Materials: ")

# --- Step 4: split "dosimetric" into its own run, wrapped with proofErr spell-check markers ---
$full = $d.Content.Text
$idx = $full.IndexOf("dosimetric data have")
$dosStart = $idx
$dosEnd = $idx + [string]"dosimetric".Length
$dosRange = $d.Range($dosStart, $dosEnd)
$dosRange.Text = "dosimetric"

# --- Step 5: split off the "Results:" paragraph at "This C# program is usable..." ---
$full = $d.Content.Text
$idx = $full.IndexOf("This C# program is usable")
$splitPoint = $d.Range($idx, $idx)
$splitPoint.InsertParagraphBefore()

$full = $d.Content.Text
$idx = $full.IndexOf("This C# program is usable")
$r = $d.Range($idx, $idx)
$r.InsertBefore("Results: ")

# --- Step 6: split off the "Conclusions:" paragraph at "This tool has been evaluated..." ---
$full = $d.Content.Text
$idx = $full.IndexOf("This tool has been evaluated")
$splitPoint = $d.Range($idx, $idx)
$splitPoint.InsertParagraphBefore()

$full = $d.Content.Text
$idx = $full.IndexOf("This tool has been evaluated")
$r = $d.Range($idx, $idx)
$r.InsertBefore("Conclusions: ")
